$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.695.15"
$ws.Range("E2").Value = "'  +0.56%  "
$ws.Range("D3").Value = "'1.850.86"
$ws.Range("E3").Value = "'  +0.50%  "
$ws.Range("E4").Value = "'  +0.36%  "
$ws.Range("D5").Value = "'312.72"
$ws.Range("E5").Value = "'  +0.02%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.44%  "
$ws.Range("D7").Value = "'0.4270"
$ws.Range("E7").Value = "'  +0.55%  "
$ws.Range("E8").Value = "'  +0.28%  "
$ws.Range("D9").Value = "'44.80"
$ws.Range("E9").Value = "'  +3.03%  "
$ws.Range("E10").Value = "'  +1.53%  "
$ws.Range("D11").Value = "'0.8769"
$ws.Range("E11").Value = "'  -2.29%  "
$ws.Range("D12").Value = "'20.65"
$ws.Range("E12").Value = "'  -0.12%  "
$ws.Range("D13").Value = "'1.877.25"
$ws.Range("E13").Value = "'  +1.59%  "
$ws.Range("D14").Value = "'5.325"
$ws.Range("E14").Value = "'  +0.17%  "
$ws.Range("D15").Value = "'6.525"
$ws.Range("E15").Value = "'  -0.90%  "
$ws.Range("D16").Value = "'0.06892"
$ws.Range("E16").Value = "'  +1.21%  "
$ws.Range("E17").Value = "'  +0.50%  "
$ws.Range("E18").Value = "'  +3.38%  "
$ws.Range("D19").Value = "'0.000009034"
$ws.Range("E19").Value = "'  +1.26%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "'  +0.48%  "
$ws.Range("E21").Value = "'  +0.00%  "
$ws.Range("D22").Value = "'27.715.56"
$ws.Range("E22").Value = "'  +0.67%  "
$ws.Range("D23").Value = "'4.967"
$ws.Range("E23").Value = "'  +0.66%  "
$ws.Range("E24").Value = "'  -3.25%  "
$ws.Range("D25").Value = "'2.157.88"
$ws.Range("E25").Value = "'  +4.94%  "
$ws.Range("E26").Value = "'  -3.80%  "
$ws.Range("D27").Value = "'154.05"
$ws.Range("E27").Value = "'  +1.66%  "
$ws.Range("D28").Value = "'18.80"
$ws.Range("E28").Value = "'  +3.46%  "
$ws.Range("D29").Value = "'121.52"
$ws.Range("E29").Value = "'  +9.43%  "
$ws.Range("D30").Value = "'5.267"
$ws.Range("E30").Value = "'  -0.90%  "
$ws.Range("E31").Value = "'  +7.76%  "
$ws.Range("D32").Value = "'0.08917"
$ws.Range("E32").Value = "'  +0.40%  "
$ws.Range("D33").Value = "'0.7608"
$ws.Range("E33").Value = "'  -1.92%  "
$ws.Range("D34").Value = "'2.969"
$ws.Range("E34").Value = "'  +4.13%  "
$ws.Range("D35").Value = "'4.523"
$ws.Range("E35").Value = "'  +1.20%  "
$ws.Range("D36").Value = "'1.102"
$ws.Range("E36").Value = "'  +1.75%  "
$ws.Range("D37").Value = "'0.05393"
$ws.Range("E37").Value = "'  -0.03%  "
$ws.Range("E38").Value = "'  -0.41%  "
$ws.Range("D39").Value = "'0.01934"
$ws.Range("E39").Value = "'  +0.78%  "
$ws.Range("D40").Value = "'2.820"
$ws.Range("E40").Value = "'  -4.33%  "
$ws.Range("D41").Value = "'0.5080"
$ws.Range("E41").Value = "'  +0.77%  "
$ws.Range("E42").Value = "'  +1.19%  "
$ws.Range("D43").Value = "'6.779"
$ws.Range("E43").Value = "'  -0.12%  "
$ws.Range("D44").Value = "'8.346"
$ws.Range("E44").Value = "'  +1.59%  "
$ws.Range("D45").Value = "'0.06543"
$ws.Range("E45").Value = "'  -0.90%  "
$ws.Range("D46").Value = "'10.34"
$ws.Range("B47").Value = "'Decentraland"
$ws.Range("C47").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4684"
$ws.Range("E47").Value = "'  -0.45%  "
$ws.Range("B48").Value = "'Quant"
$ws.Range("C48").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'105.06"
$ws.Range("E48").Value = "'  -1.00%  "
$ws.Range("E49").Value = "'  +0.53%  "
$ws.Range("E50").Value = "'  -1.31%  "
$ws.Range("B51").Value = "'Aave"
$ws.Range("C51").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'64.59"
$ws.Range("E51").Value = "'  +0.13%  "
